$d = $word.ActiveDocument

# Locate the target paragraph (bullet under "Creative Team Contributions:")
# that currently reads:
#   "Who are the most frequent directors, writers, and composers in Pixar's history?"
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Who are the most frequent directors*") {
        $targetPara = $p
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not locate the target paragraph."
}

$pRange = $targetPara.Range

# Pull the paragraph's real OOXML so we can keep its opening <w:p ...>
# attributes and its <w:pPr> (list/numbering formatting etc.) completely
# untouched - only the runs inside the paragraph should change.
$packageXml = $pRange.WordOpenXML
if ($packageXml -notmatch '(?s)<w:p\b.*?</w:p>') {
    throw "Could not extract paragraph XML."
}
$paraXml = $matches[0]

if ($paraXml -match '(?s)^(.*</w:pPr>)(.*)</w:p>$') {
    $prefix = $matches[1]
} elseif ($paraXml -match '(?s)^(<w:p\b[^>]*>)(.*)</w:p>$') {
    $prefix = $matches[1]
} else {
    throw "Could not split paragraph XML into pPr/runs."
}

# Replace the sentence with the same wording split across five runs:
# "...directors, writers, and composers in Pixar's history?"
#   -> "...directors and writers in Pixar's history?"
$newRuns =
    '<w:r><w:t>Who are the most frequent directors</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
    '<w:r><w:t>writers</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>in Pixar''s history?</w:t></w:r>'

$newParaXml = $prefix + $newRuns + '</w:p>'

$pRange.InsertXML($newParaXml)
